$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sampling schemes added by the averaged-intensities run for the
# spiral sampling schemes. Append them as three new rows (17-19),
# mirroring the existing table layout (index column A, scheme-name
# column B, then the per-HKL averaged-intensity values in C:M, all 1).

$newSchemes = @(
    @{ Row = 17; Idx = 15; Name = "Spiral-90deg-10rot-5space" },
    @{ Row = 18; Idx = 16; Name = "Spiral-90deg-15rot-5space" },
    @{ Row = 19; Idx = 17; Name = "Spiral-90deg-10rot-3space" }
)

foreach ($scheme in $newSchemes) {
    $r = $scheme.Row

    # Clone formatting (bold/centered/bordered index style + plain data
    # cells) from the last existing data row so the new rows look exactly
    # like the others, then overwrite the values.
    $ws.Range("A16:M16").Copy($ws.Range("A" + $r + ":M" + $r))

    $ws.Range("A$r").Value = $scheme.Idx
    $ws.Range("B$r").Value = $scheme.Name

    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($r, $col).Value = 1
    }
}
